# Updated cryptos list -- apply latest Price (D) and Volume(1h) (E) figures,
# plus the RenderToken/Stacks row swap (rows 37-38), per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to remain a text value (matches the workbook's existing
    # inline-string cells) even when the new value looks numeric (e.g. "546.77"),
    # then drop the temporary text format so no stray style is left behind.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "67.086.66"
$ws.Range("E2").Value = "  -3.07%  "

# Row 3
Set-TextValue "D3" "2.385.78"
$ws.Range("E3").Value = "  -3.32%  "

# Row 5
Set-TextValue "D5" "546.77"
$ws.Range("E5").Value = "  -2.13%  "

# Row 6
Set-TextValue "D6" "156.24"
$ws.Range("E6").Value = "  -4.28%  "

# Row 8
Set-TextValue "D8" "0.501"
$ws.Range("E8").Value = "  -0.49%  "

# Row 9
Set-TextValue "D9" "0.154"
$ws.Range("E9").Value = "  +0.88%  "

# Row 10
$ws.Range("E10").Value = "  -1.26%  "

# Row 11
Set-TextValue "D11" "0.323"
$ws.Range("E11").Value = "  -3.54%  "

# Row 12
Set-TextValue "D12" "4.72"
$ws.Range("E12").Value = "  -2.14%  "

# Row 13
Set-TextValue "D13" "67.023.06"
$ws.Range("E13").Value = "  -3.04%  "

# Row 14
Set-TextValue "D14" "0.0000166"
$ws.Range("E14").Value = "  -2.11%  "

# Row 15
Set-TextValue "D15" "22.62"
$ws.Range("E15").Value = "  -4.26%  "

# Row 16
Set-TextValue "D16" "10.22"
$ws.Range("E16").Value = "  -5.08%  "

# Row 17
Set-TextValue "D17" "325.52"
$ws.Range("E17").Value = "  -5.07%  "

# Row 18
Set-TextValue "D18" "6.72"
$ws.Range("E18").Value = "  -4.87%  "

# Row 19
Set-TextValue "D19" "3.72"
$ws.Range("E19").Value = "  -2.04%  "

# Row 20
$ws.Range("E20").Value = "  -0.51%  "

# Row 21
Set-TextValue "D21" "1.82"
$ws.Range("E21").Value = "  -5.46%  "

# Row 22
Set-TextValue "D22" "65.31"
$ws.Range("E22").Value = "  -2.70%  "

# Row 23
Set-TextValue "D23" "3.55"
$ws.Range("E23").Value = "  -3.56%  "

# Row 24
Set-TextValue "D24" "7.90"
$ws.Range("E24").Value = "  -3.61%  "

# Row 25
Set-TextValue "D25" "0.0$([char]0x2083)0785"
$ws.Range("E25").Value = "  -4.22%  "

# Row 26
Set-TextValue "D26" "6.91"
$ws.Range("E26").Value = "  -3.77%  "

# Row 27
$ws.Range("E27").Value = "  +0.21%  "

# Row 28
Set-TextValue "D28" "411.03"
$ws.Range("E28").Value = "  -6.47%  "

# Row 29
Set-TextValue "D29" "1.10"
$ws.Range("E29").Value = "  -3.33%  "

# Row 30
Set-TextValue "D30" "1.57"
$ws.Range("E30").Value = "  -2.92%  "

# Row 31
Set-TextValue "D31" "159.21"
$ws.Range("E31").Value = "  +1.58%  "

# Row 32
$ws.Range("E32").Value = "  -0.76%  "

# Row 33
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
Set-TextValue "D34" "17.56"
$ws.Range("E34").Value = "  -1.89%  "

# Row 35
$ws.Range("E35").Value = "  -5.26%  "

# Row 36
Set-TextValue "D36" "0.289"
$ws.Range("E36").Value = "  -4.45%  "

# Row 37
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D37" "1.43"
$ws.Range("E37").Value = "  -3.25%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D38" "4.15"
$ws.Range("E38").Value = "  -6.80%  "

# Row 39
$ws.Range("E39").Value = "  -4.42%  "

# Row 40
Set-TextValue "D40" "1.94"
$ws.Range("E40").Value = "  -6.72%  "

# Row 41
Set-TextValue "D41" "3.25"
$ws.Range("E41").Value = "  -3.22%  "

# Row 42
Set-TextValue "D42" "127.37"
$ws.Range("E42").Value = "  -4.20%  "

# Row 43
Set-TextValue "D43" "0.0702"
$ws.Range("E43").Value = "  -2.37%  "

# Row 44
Set-TextValue "D44" "0.467"
$ws.Range("E44").Value = "  -3.35%  "

# Row 45
$ws.Range("E45").Value = "  -2.84%  "

# Row 46
Set-TextValue "D46" "0.0905"
$ws.Range("E46").Value = "  -0.22%  "

# Row 47
Set-TextValue "D47" "1.10"
$ws.Range("E47").Value = "  -1.12%  "

# Row 48
Set-TextValue "D48" "1.31"
$ws.Range("E48").Value = "  -8.74%  "

# Row 49
Set-TextValue "D49" "16.21"
$ws.Range("E49").Value = "  -4.11%  "

# Row 50
Set-TextValue "D50" "0.0$([char]0x2086)0200"
$ws.Range("E50").Value = "  -4.85%  "

# Row 51
Set-TextValue "D51" "0.0420"
$ws.Range("E51").Value = "  -2.40%  "
